# fwt_case5.xlsx — adjust test fixture for _make_route:
#   - fix the posting date on row 3 (was a duplicate of row 2's date)
#   - autofit column A ("Posting Date") now that its header is the
#     longest value in the column
#   - leave the cursor where the user last clicked while reviewing the
#     sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3's "Posting Date" (A3) incorrectly duplicated row 2's date (44950 /
# 2023-01-24). Bump it to the correct next day (44951 / 2023-01-25).
$ws.Range("A3").Value = 44951

# Column A now needs to be wide enough to show its "Posting Date" header
# in full, so autofit it to its contents (mirrors bestFit width behavior).
$ws.Columns("A").AutoFit() | Out-Null

# Restore the last active selection on the sheet.
$ws.Range("Q11").Select() | Out-Null
